$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F, copying header style from E1
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Populate time_taken values for each data row
$ws.Cells.Item(2, 6).Value = "2021-10-05 13:41:47.607655"
$ws.Cells.Item(3, 6).Value = "2021-10-05 13:41:47.607667"
$ws.Cells.Item(4, 6).Value = "2021-10-05 13:41:47.607672"
$ws.Cells.Item(5, 6).Value = "2021-10-05 13:41:47.607675"
$ws.Cells.Item(6, 6).Value = "2021-10-05 13:41:47.607678"
$ws.Cells.Item(7, 6).Value = "2021-10-05 13:41:47.607681"
$ws.Cells.Item(8, 6).Value = "2021-10-05 13:41:47.607684"
$ws.Cells.Item(9, 6).Value = "2021-10-05 13:41:47.607687"
$ws.Cells.Item(10, 6).Value = "2021-10-05 13:41:47.607690"
$ws.Cells.Item(11, 6).Value = "2021-10-05 13:41:47.607693"
$ws.Cells.Item(12, 6).Value = "2021-10-05 13:41:47.607696"
$ws.Cells.Item(13, 6).Value = "2021-10-05 13:41:47.607699"
$ws.Cells.Item(14, 6).Value = "2021-10-05 13:41:47.607702"
$ws.Cells.Item(15, 6).Value = "2021-10-05 13:41:47.607705"
$ws.Cells.Item(16, 6).Value = "2021-10-05 13:41:47.607709"
$ws.Cells.Item(17, 6).Value = "2021-10-05 13:41:47.607712"
$ws.Cells.Item(18, 6).Value = "2021-10-05 13:41:47.607715"
$ws.Cells.Item(19, 6).Value = "2021-10-05 13:41:47.607718"
$ws.Cells.Item(20, 6).Value = "2021-10-05 13:41:47.607721"
$ws.Cells.Item(21, 6).Value = "2021-10-05 13:41:47.607725"
$ws.Cells.Item(22, 6).Value = "2021-10-05 13:41:47.607728"

$excel.CutCopyMode = 0
